# login node Stateless service
#
# Inserts 6 new login-error tip rows (TooManyDevices, LoginFSMLoadFailed,
# LoginFSMEventFailed, LoginAccountDataLoadFaile, LoginSessionNotFound,
# LoginAccountDataLoadFailed) right after the existing "LoginRedisSetFailed"
# row and before the "//scene_error" section header. Inserting whole rows
# pushes every following row down by 6, which is the only structural change
# needed - all later rows keep their original cell content/relationships.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
  "    TooManyDevices",
  "    LoginFSMLoadFailed",
  "    LoginFSMEventFailed",
  "    LoginAccountDataLoadFaile",
  "    LoginSessionNotFound",
  "    LoginAccountDataLoadFailed"
)

$firstRow = 63

# Insert 6 blank rows above the current row 63 ("//scene_error"), shifting
# it (and everything below) down to row 69.
$ws.Rows.Item($firstRow).Resize($newValues.Length).Insert()

for ($i = 0; $i -lt $newValues.Length; $i++) {
  $row = $firstRow + $i
  $ws.Cells.Item($row, 1).Value = $newValues[$i]
}

# Leave the cursor on the last inserted row, matching the author's final
# selection after making the edit.
$ws.Range("A68").Select()
